$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet to reflect the new revision (v1.3 -> v1.3.1)
$ws.Name = "Mini HUB v1.3.1"

# 2. Update the title in A1 (rich text: bold run + italic run).
#    Only the first (bold) run's text changes, "v1.3" -> "v1.3.1".
$titleCell = $ws.Range("A1")
$titleCell.Characters(1, 43).Text = "Chordata Mini HUB Bill of Materials - v1.3.1`n"
$titleCell.Characters(1, 45).Font.Bold = $true
$titleCell.Characters(46, 10).Font.Italic = $true

# 3. Remove the two unused 0805 MLCC capacitor lines (items 8 & 9, rows 12-13).
#    Clear all their data and mark the designator column as "<removed>"
#    in an italic, left-aligned style.
$ws.Range("B12:N12").ClearContents()
$ws.Range("B13:N13").ClearContents()

$row12B = $ws.Range("B12")
$row12B.Value = "<removed>"
$row12B.Font.Italic = $true
$row12B.HorizontalAlignment = -4131

$row13B = $ws.Range("B13")
$row13B.Value = "<removed>"
$row13B.Font.Italic = $true
$row13B.HorizontalAlignment = -4131

# 4. Move the selection to the merged "Item" header cell (A2:A4).
$ws.Range("A2:A4").Select()
